# Swap the contents of columns C and D (codeforiati:group-code /
# codeforiati:group-name) for every row in the sheet, including the header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

$colC = $ws.Range("C1:C$lastRow").Value2
$colD = $ws.Range("D1:D$lastRow").Value2

$ws.Range("C1:C$lastRow").Value2 = $colD
$ws.Range("D1:D$lastRow").Value2 = $colC
